$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$ws.Range("A2").Value = 20330051920252
$ws.Range("B2").Value = "SANCHEZ"
$ws.Range("C2").Value = "PEREZ"
$ws.Range("D2").Value = "ARLET"
$ws.Range("E2").Value = "INGLÉS III"
$ws.Range("F2").Value = "3ALCM"
$ws.Range("G2").Value = 6

$ws.Range("A3").Value = 20330051920263
$ws.Range("B3").Value = "CARRERA"
$ws.Range("C3").Value = "ZAVALETA"
$ws.Range("D3").Value = "ALFREDO"
$ws.Range("E3").Value = "INGLÉS III"
$ws.Range("F3").Value = "3APM"
$ws.Range("G3").Value = 6

$ws.Range("A4").Value = 20330051920040
$ws.Range("B4").Value = "AMADOR"
$ws.Range("C4").Value = "PORRAS"
$ws.Range("D4").Value = "FRANCISCO ALAN"
$ws.Range("E4").Value = "INGLÉS III"
$ws.Range("F4").Value = "3BEM"
$ws.Range("G4").Value = 6

$ws.Range("A5").Value = 20330051920284
$ws.Range("B5").Value = "ANTONIO"
$ws.Range("C5").Value = "GARCIA"
$ws.Range("D5").Value = "ISRAEL"
$ws.Range("E5").Value = "INGLÉS III"
$ws.Range("F5").Value = "3BLCM"
$ws.Range("G5").Value = 6

$ws.Range("A6").Value = 20330051920295
$ws.Range("B6").Value = "HERRERA"
$ws.Range("C6").Value = "CERON"
$ws.Range("D6").Value = "YAMILE"
$ws.Range("E6").Value = "INGLÉS III"
$ws.Range("F6").Value = "3BLCM"
$ws.Range("G6").Value = 6

$ws.Range("A7").Value = 20330051920313
$ws.Range("B7").Value = "TORRES"
$ws.Range("C7").Value = "VAZQUEZ"
$ws.Range("D7").Value = "JOSELIN GUADALUPE"
$ws.Range("E7").Value = "INGLÉS III"
$ws.Range("F7").Value = "3BLCM"
$ws.Range("G7").Value = 6

$ws.Range("A8").Value = 20330051920373
$ws.Range("B8").Value = "RICO"
$ws.Range("C8").Value = "BAUTISTA"
$ws.Range("D8").Value = "EDGAR RAMSES"
$ws.Range("E8").Value = "INGLÉS III"
$ws.Range("F8").Value = "3ARHM"
$ws.Range("G8").Value = 6

$ws.Range("A9").Value = 20330051920053
$ws.Range("B9").Value = "GARCIA"
$ws.Range("C9").Value = "MARTINEZ"
$ws.Range("D9").Value = "MARIA ASUNCION"
$ws.Range("E9").Value = "INGLÉS III"
$ws.Range("F9").Value = "3BEM"
$ws.Range("G9").Value = 6

